# Apply updated crypto price/volume data to the active worksheet
# (values refreshed by the scheduled GitHub Actions scrape run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.273.44'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '1.861.48'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7043'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08197'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3051'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08182'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").Value = '1.836.32'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7181'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.191'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '29.283.79'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.786'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.17%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007859'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9990'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = '2.105.03'
$ws.Range("E22").Value = '  +2.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.473'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.010'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1457'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.437'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.45%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.420'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.483'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05223'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.173'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7093'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9989'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.663'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01850'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.723'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.84%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9269'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.96%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.145.82'
$ws.Range("E42").Value = '  +8.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.955'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4280'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.774'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("D49").Value = '2.001.66'
$ws.Range("E49").Value = '  +1.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.205'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.972'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.75%  '
